$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 800.1111
$ws.Range("I9").Value = 800.1111
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 800.1111
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -631.1111
$ws.Range("N9").Value = $null
$ws.Range("H16").Value = 4252.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 4252.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 4252.5
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = -4712.5
$ws.Range("H18").Value = 1814.9
$ws.Range("I18").Value = 1381.125
$ws.Range("J18").Value = 3550
$ws.Range("K18").Value = 1381.125
$ws.Range("L18").Value = 3550
$ws.Range("M18").Value = -1097.125
$ws.Range("N18").Value = -4118
$ws.Range("H51").Value = 2
$ws.Range("J51").Value = 2
$ws.Range("L51").Value = 2
$ws.Range("N51").Value = -970
$ws.Range("H62").Value = 5886.8887
$ws.Range("I62").Value = 2992
$ws.Range("K62").Value = 2992
$ws.Range("M62").Value = -2368
$ws.Range("H65").Value = 5886.8887
$ws.Range("I65").Value = 2992
$ws.Range("K65").Value = 14960
$ws.Range("M65").Value = -11840
$ws.Range("H126").Value = 70439.86
$ws.Range("J126").Value = 67019.75
$ws.Range("L126").Value = 67019.75
$ws.Range("N126").Value = -76899.75
$ws.Range("H132").Value = 758.125
$ws.Range("I132").Value = 671.7895
$ws.Range("K132").Value = 2015.3685
$ws.Range("M132").Value = 514.6315
$ws.Range("H137").Value = 3690.868
$ws.Range("I137").Value = 2728.147
$ws.Range("J137").Value = 5413.6313
$ws.Range("K137").Value = 8184.441
$ws.Range("L137").Value = 16240.8939
$ws.Range("M137").Value = -5634.441
$ws.Range("N137").Value = -21340.8939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8776563
$ws.Range("I32").Value = 10871615
$ws.Range("K32").Value = 10871615
$ws.Range("M32").Value = -10871328
$ws.Range("H63").Value = 7951.1
$ws.Range("I63").Value = 4499.5
$ws.Range("K63").Value = 4499.5
$ws.Range("M63").Value = -3813.5
$ws.Range("H66").Value = 7951.1
$ws.Range("I66").Value = 4499.5
$ws.Range("K66").Value = 22497.5
$ws.Range("M66").Value = -19065.5
$ws.Range("H92").Value = 55497.6
$ws.Range("J92").Value = 55497.6
$ws.Range("L92").Value = 55497.6
$ws.Range("N92").Value = -60489.6
$ws.Range("H138").Value = 182498
$ws.Range("J138").Value = 191997.6
$ws.Range("L138").Value = 191997.6
$ws.Range("N138").Value = -202277.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2557.9167
$ws.Range("I99").Value = 1831.2
$ws.Range("J99").Value = 3769.111
$ws.Range("K99").Value = 1831.2
$ws.Range("L99").Value = 3769.111
$ws.Range("M99").Value = -333.2
$ws.Range("N99").Value = -6765.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 100000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 100000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -100226
$ws.Range("H3").Value = 1250
$ws.Range("J3").Value = 1333.3334
$ws.Range("L3").Value = 1333.3334
$ws.Range("N3").Value = -1559.3334
$ws.Range("H5").Value = 11751
$ws.Range("J5").Value = 14999
$ws.Range("L5").Value = 14999
$ws.Range("N5").Value = -15223
$ws.Range("H8").Value = 8015
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H10").Value = 10262.833
$ws.Range("I10").Value = 1923.5
$ws.Range("J10").Value = 18602.166
$ws.Range("K10").Value = 1923.5
$ws.Range("L10").Value = 18602.166
$ws.Range("M10").Value = -1784.5
$ws.Range("N10").Value = -18880.166
$ws.Range("H11").Value = 44003
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 58337.332
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 58337.332
$ws.Range("M11").Value = -860
$ws.Range("N11").Value = -58617.332
$ws.Range("H13").Value = 3401
$ws.Range("I13").Value = 1552
$ws.Range("K13").Value = 1552
$ws.Range("M13").Value = -1413
$ws.Range("H14").Value = 13656.571
$ws.Range("I14").Value = 1549.5
$ws.Range("J14").Value = 18499.4
$ws.Range("K14").Value = 1549.5
$ws.Range("L14").Value = 18499.4
$ws.Range("M14").Value = -1379.5
$ws.Range("N14").Value = -18839.4
$ws.Range("H19").Value = 2547.111
$ws.Range("I19").Value = 3145
$ws.Range("J19").Value = 1799.75
$ws.Range("K19").Value = 3145
$ws.Range("L19").Value = 1799.75
$ws.Range("M19").Value = -2975
$ws.Range("N19").Value = -2139.75
$ws.Range("H24").Value = 2547.111
$ws.Range("I24").Value = 3145
$ws.Range("J24").Value = 1799.75
$ws.Range("K24").Value = 3145
$ws.Range("L24").Value = 1799.75
$ws.Range("M24").Value = -2975
$ws.Range("N24").Value = -2139.75
$ws.Range("H31").Value = 656761
$ws.Range("I31").Value = 10189.35
$ws.Range("J31").Value = 1464975.5
$ws.Range("K31").Value = 10189.35
$ws.Range("L31").Value = 1464975.5
$ws.Range("M31").Value = -9894.35
$ws.Range("N31").Value = -1465565.5
$ws.Range("H34").Value = 656761
$ws.Range("I34").Value = 10189.35
$ws.Range("J34").Value = 1464975.5
$ws.Range("K34").Value = 10189.35
$ws.Range("L34").Value = 1464975.5
$ws.Range("M34").Value = -9987.35
$ws.Range("N34").Value = -1465379.5
$ws.Range("H107").Value = 1784
$ws.Range("I107").Value = 735.4706
$ws.Range("J107").Value = 3155.1538
$ws.Range("K107").Value = 735.4706
$ws.Range("L107").Value = 3155.1538
$ws.Range("M107").Value = 1184.5294
$ws.Range("N107").Value = -6995.1538
$ws.Range("H122").Value = 2004
$ws.Range("I122").Value = 1654.3125
$ws.Range("K122").Value = 4962.9375
$ws.Range("M122").Value = -2512.9375
$ws.Range("H133").Value = 71333.336
$ws.Range("J133").Value = 71333.336
$ws.Range("L133").Value = 71333.336
$ws.Range("N133").Value = -76393.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 218.58333
$ws.Range("I15").Value = 115.375
$ws.Range("K15").Value = 346.125
$ws.Range("M15").Value = -206.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null
$ws.Range("H107").Value = 745.86664
$ws.Range("I107").Value = 730.125
$ws.Range("J107").Value = 763.8570999999999
$ws.Range("K107").Value = 730.125
$ws.Range("L107").Value = 763.8570999999999
$ws.Range("M107").Value = 1189.875
$ws.Range("N107").Value = -4603.8571
$ws.Range("H126").Value = 4472.727
$ws.Range("I126").Value = 4070.6667
$ws.Range("J126").Value = 4623.5
$ws.Range("K126").Value = 12212.0001
$ws.Range("L126").Value = 13870.5
$ws.Range("M126").Value = -9742.000100000001
$ws.Range("N126").Value = -18810.5
$ws.Range("H132").Value = 16671639
$ws.Range("I132").Value = 25002726
$ws.Range("J132").Value = 9465.15
$ws.Range("K132").Value = 75008178
$ws.Range("L132").Value = 28395.45
$ws.Range("M132").Value = -75005648
$ws.Range("N132").Value = -33455.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 151003.42
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 151003.42
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 151003.42
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -151227.42
$ws.Range("H40").Value = 3525.52
$ws.Range("I40").Value = 2448.2856
$ws.Range("J40").Value = 3944.4443
$ws.Range("K40").Value = 2448.2856
$ws.Range("L40").Value = 3944.4443
$ws.Range("M40").Value = -2312.2856
$ws.Range("N40").Value = -4216.4443
$ws.Range("H68").Value = 3000.6667
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("H71").Value = 3000.6667
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
$ws.Range("H109").Value = 101640
$ws.Range("J109").Value = 101640
$ws.Range("L109").Value = 101640
$ws.Range("N109").Value = -104414
$ws.Range("H126").Value = 151003.42
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 151003.42
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 453010.26
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -457950.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 23500
$ws.Range("J11").Value = 25000
$ws.Range("L11").Value = 25000
$ws.Range("N11").Value = -25284
$ws.Range("H24").Value = 83343.336
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H100").Value = 585.72
$ws.Range("I100").Value = 405.66666
$ws.Range("K100").Value = 811.33332
$ws.Range("M100").Value = -270.33332
$ws.Range("H126").Value = 1249.909
$ws.Range("I126").Value = 1224.9
$ws.Range("K126").Value = 3674.7
$ws.Range("M126").Value = -1204.7
